# "fix model, init data"
#
# The "Node" column (header "Node" in E6, data placeholder
# "[[Data.Node.Name]]" in E7) is removed from the BusStop export template.
# The "DeliveryTrips" column (header + "[[Data.DeliveryTrips]]" placeholder),
# which sat one column to the right (F), shifts left into E so the table has
# no gap.
#
# Concretely, on Sheet1:
#   Row 6 (headers):      E6 <- F6 ("DeliveryTrips"), F6 removed
#   Row 7 (placeholders): E7 <- F7 ("[[Data.DeliveryTrips]]"), F7 removed
#
# Rows 1, 3, 4 and 5 (logo, title block, spacer) are untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Move "DeliveryTrips" header (was F6) into E6, overwriting the old "Node"
# header that used to live there.
$ws.Range("E6").Value2 = $ws.Range("F6").Value2

# Move the "[[Data.DeliveryTrips]]" placeholder (was F7) into E7, overwriting
# the old "[[Data.Node.Name]]" placeholder.
$ws.Range("E7").Value2 = $ws.Range("F7").Value2

# Remove the now-duplicated trailing cells so the row ends at column E again.
$ws.Range("F6").Clear()
$ws.Range("F7").Clear()
